# Katalon AI generated update:
#  - Renamed data-binding placeholder labels in row 1 (div_testRunDetails_* -> div_testRuns_*,
#    link_projectLinks_* -> link_testProject_*)
#  - Updated the referenced test-data file path in A2
#    (scheduleAndRunTestSuite-test-data -> scheduleAndRunTestWithConfigurations-test-data)
#  - Column widths are re-fit to the new (shorter/longer) header & data text, matching
#    Excel's AutoFit behaviour (stored column width = ColumnWidth + 0.8333333...,
#    so we subtract that padding to land on the exact target integer widths).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell text -------------------------------------------------

$ws.Range("B1").Value = "div_testRuns_internalRoleCellName"

$ws.Range("D1").Value = "link_testProject_internalRoleLinkName"
$ws.Range("E1").Value = "link_testProject_internalRoleLinkName_1"
$ws.Range("F1").Value = "link_testProject_project_id"
$ws.Range("G1").Value = "link_testProject_project_id_1"
$ws.Range("H1").Value = "link_testProject_team_id"
$ws.Range("I1").Value = "link_testProject_team_id_1"
$ws.Range("J1").Value = "link_testProject_test_project_id"
$ws.Range("K1").Value = "link_testProject_test_project_id_1"
$ws.Range("L1").Value = "link_testProject_trNthChild"
$ws.Range("M1").Value = "link_testProject_trNthChild_1"

$ws.Range("A2").Value = "Data Files/AI-Generated/Common/scheduleAndRunTestWithConfigurations-test-data"

# --- Re-fit column widths to match the new text ------------------------
# Setting .ColumnWidth via COM stores (ColumnWidth + 5/6) as the sheet's
# <col width> attribute, so subtract the padding to hit the exact target.

$padding = 5 / 6

$targetWidths = @{
    1  = 79
    2  = 35
    4  = 39
    5  = 41
    6  = 29
    7  = 31
    8  = 26
    9  = 28
    10 = 34
    11 = 36
    12 = 29
    13 = 31
}

foreach ($col in $targetWidths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $targetWidths[$col] - $padding
}
